# Update derived data map xlsx
#
# The "sheet1" worksheet tracked one row per (study, breakdown) pair that
# feeds into the derived-data pipeline. Several studies/rows that are no
# longer part of the pipeline are removed (CHE1 "region" breakdown row,
# the NYC_NY_1 row, the BRA5 row, the LA_CA1 row, and the NYC_NY_1_nch
# care-home row), and the DNK1 row is switched from the "region" breakdown
# to the "ageband" breakdown (updating its relpath to match).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet1")

# Delete obsolete rows. Deleting from the bottom up keeps the remaining
# row numbers stable as each EntireRow.Delete() shifts rows below it up.
$ws.Rows.Item(23).EntireRow.Delete()   # NYC_NY_1_nch (care_home_deaths = no)
$ws.Rows.Item(19).EntireRow.Delete()   # LA_CA1 / region
$ws.Rows.Item(18).EntireRow.Delete()   # BRA5 / region
$ws.Rows.Item(10).EntireRow.Delete()   # NYC_NY_1 / ageband
$ws.Rows.Item(7).EntireRow.Delete()    # CHE1 / region

# DNK1 (now row 5) switches breakdown from "region" to "ageband".
$ws.Range("C5").Value = "ageband"
$ws.Range("D5").Value = "data/derived/DNK1/DNK1_agebands.RDS"

# Move the active selection, matching the saved cursor position.
[void]$ws.Range("A9").Select()
